# Update crypto price/volume data (scraped refresh) across rows 2-51
# Rows 48/49 also swap Coin/Link (Aave <-> Frax) per source reorder.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "26.007.01"
$cell.ClearFormats()
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.54%  "
$cell.ClearFormats()

# Row 3
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.639.40"
$cell.ClearFormats()
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.19%  "
$cell.ClearFormats()

# Row 4
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.22%  "
$cell.ClearFormats()

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "214.78"
$cell.ClearFormats()
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.15%  "
$cell.ClearFormats()

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.5087"
$cell.ClearFormats()
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.16%  "
$cell.ClearFormats()

# Row 7
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.01%  "
$cell.ClearFormats()

# Row 8
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.2564"
$cell.ClearFormats()
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.08%  "
$cell.ClearFormats()

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.06369"
$cell.ClearFormats()
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.05%  "
$cell.ClearFormats()

# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "19.50"
$cell.ClearFormats()
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.05%  "
$cell.ClearFormats()

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.07756"
$cell.ClearFormats()
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.22%  "
$cell.ClearFormats()

# Row 12
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.667.58"
$cell.ClearFormats()
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.70%  "
$cell.ClearFormats()

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.286"
$cell.ClearFormats()
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.58%  "
$cell.ClearFormats()

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.5443"
$cell.ClearFormats()
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.55%  "
$cell.ClearFormats()

# Row 15
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0₅7735"
$cell.ClearFormats()
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.65%  "
$cell.ClearFormats()

# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "64.21"
$cell.ClearFormats()
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.56%  "
$cell.ClearFormats()

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "26.047.19"
$cell.ClearFormats()
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.58%  "
$cell.ClearFormats()

# Row 18
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.01%  "
$cell.ClearFormats()

# Row 19
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "196.98"
$cell.ClearFormats()
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.40%  "
$cell.ClearFormats()

# Row 20
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.424"
$cell.ClearFormats()
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.09%  "
$cell.ClearFormats()

# Row 21
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.14%  "
$cell.ClearFormats()

# Row 22
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.037"
$cell.ClearFormats()
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.29%  "
$cell.ClearFormats()

# Row 23
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.03%  "
$cell.ClearFormats()

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.879"
$cell.ClearFormats()
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.23%  "
$cell.ClearFormats()

# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "141.23"
$cell.ClearFormats()
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.92%  "
$cell.ClearFormats()

# Row 26
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.07%  "
$cell.ClearFormats()

# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.837"
$cell.ClearFormats()
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.04%  "
$cell.ClearFormats()

# Row 28
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "15.52"
$cell.ClearFormats()
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.85%  "
$cell.ClearFormats()

# Row 29
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.13%  "
$cell.ClearFormats()

# Row 30
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.04867"
$cell.ClearFormats()
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.03%  "
$cell.ClearFormats()

# Row 31
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.254"
$cell.ClearFormats()
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.03%  "
$cell.ClearFormats()

# Row 32
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.167"
$cell.ClearFormats()
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.53%  "
$cell.ClearFormats()

# Row 33
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.16%  "
$cell.ClearFormats()

# Row 34
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.362"
$cell.ClearFormats()
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.03%  "
$cell.ClearFormats()

# Row 35
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.8941"
$cell.ClearFormats()
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.78%  "
$cell.ClearFormats()

# Row 36
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.141.90"
$cell.ClearFormats()
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.44%  "
$cell.ClearFormats()

# Row 37
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.574"
$cell.ClearFormats()
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.82%  "
$cell.ClearFormats()

# Row 38
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.38%  "
$cell.ClearFormats()

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.01559"
$cell.ClearFormats()
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.12%  "
$cell.ClearFormats()

# Row 40
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.14%  "
$cell.ClearFormats()

# Row 41
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.518"
$cell.ClearFormats()
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.51%  "
$cell.ClearFormats()

# Row 42
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0₈127"
$cell.ClearFormats()
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +4.89%  "
$cell.ClearFormats()

# Row 43
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.8089"
$cell.ClearFormats()
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.62%  "
$cell.ClearFormats()

# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "99.15"
$cell.ClearFormats()
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.22%  "
$cell.ClearFormats()

# Row 45
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.415"
$cell.ClearFormats()
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -4.53%  "
$cell.ClearFormats()

# Row 46
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.777.53"
$cell.ClearFormats()
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.30%  "
$cell.ClearFormats()

# Row 47
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.4525"
$cell.ClearFormats()
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.56%  "
$cell.ClearFormats()

# Row 48
$cell = $ws.Cells.Item(48, 2)
$cell.NumberFormat = "@"
$cell.Value = "Frax"
$cell.ClearFormats()
$cell = $ws.Cells.Item(48, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$cell.ClearFormats()
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.ClearFormats()
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.55%  "
$cell.ClearFormats()

# Row 49
$cell = $ws.Cells.Item(49, 2)
$cell.NumberFormat = "@"
$cell.Value = "Aave"
$cell.ClearFormats()
$cell = $ws.Cells.Item(49, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell.ClearFormats()
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "54.95"
$cell.ClearFormats()
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.54%  "
$cell.ClearFormats()

# Row 50
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.05055"
$cell.ClearFormats()
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.18%  "
$cell.ClearFormats()

# Row 51
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.ClearFormats()
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.06%  "
$cell.ClearFormats()

